$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the slightly adjusted timestamp value in row 60 (A60)
$ws.Range("A60").Value = 44373.76802768287

# Add the new row 61 of data
$ws.Range("A61").Value = 44374.76922146973
$ws.Range("B61").Value = 78367
$ws.Range("C61").Value = 65845
$ws.Range("D61").Value = 3702
$ws.Range("E61").Value = 2148
$ws.Range("F61").Value = 1525
$ws.Range("G61").Value = 20756
$ws.Range("H61").Value = 1597
$ws.Range("I61").Value = 886
$ws.Range("J61").Value = 206

# Match the date-formatted number format used by the rest of column A
$ws.Range("A61").NumberFormat = $ws.Range("A60").NumberFormat
